$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (stored width = ColumnWidth + 5/6, so subtract that
# offset here to land exactly on the target stored widths of 66 and 41)
$ws.Columns.Item(1).ColumnWidth = 66 - 5/6
$ws.Columns.Item(2).ColumnWidth = 41 - 5/6

# Update cell text values
$ws.Range("B1").Value = "div_testRunDetails_internalRoleCellName"
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestSuite-test-data"
